$d = $word.ActiveDocument

# The [D10] paragraph currently ends with:
#   "...arrival/departure time."
# It needs to become:
#   "...arrival/departure time and QR code."
# i.e. insert " and QR code" right before the trailing period.
$d.Content.Find.Execute(
    "arrival/departure time.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "arrival/departure time and QR code.",
    2
) | Out-Null
